$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 39

$ws.Cells.Item($row, 1).Value = 1751
$ws.Cells.Item($row, 2).Value = "Maximum Number of Events That Can Be Attended 2"
$ws.Cells.Item($row, 3).Value = "Binary Search/Dynamic Programming"
$ws.Cells.Item($row, 4).Value = "Sort the events by ENDday, have an array of EndDays, store best result in a multi-dim array."

$ws.Range("D$row").Select()
